$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top; existing rows 1-259 shift down to 2-260
$ws.Rows.Item(1).Insert()

# Set header values
$ws.Range("A1").Value = "code"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "product_name"

# Make header row bold
$ws.Range("A1:C1").Font.Bold = $true

# Match the resulting selection seen in the diff
$ws.Range("C14").Select()
